# Updates cryptos list figures (price + 1h volume change) per the
# Fri Dec 29 20:28:08 UTC 2023 GitHub Actions data refresh.
# Column D (Price) values are numeric-looking text (e.g. "1.00", "42.049.42")
# that must stay stored as text, so they are entered with a leading
# apostrophe (the same trick Excel itself uses) to avoid Excel
# auto-converting them to numbers and dropping formatting like trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.049.42"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").Value = "'2.307.73"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'313.14"
$ws.Range("E5").Value = "  -4.02%  "
$ws.Range("D6").Value = "'105.90"
$ws.Range("E7").Value = "  -1.18%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'0.609"
$ws.Range("E9").Value = "  -2.22%  "
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("D12").Value = "'8.28"
$ws.Range("E12").Value = "  -1.66%  "
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").Value = "'0.980"
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("D15").Value = "'15.58"
$ws.Range("E15").Value = "  -5.71%  "
$ws.Range("D16").Value = "'2.658.31"
$ws.Range("E16").Value = "  -2.01%  "
$ws.Range("D17").Value = "'2.300.54"
$ws.Range("E17").Value = "  -2.04%  "
$ws.Range("D18").Value = "'42.161.96"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").Value = "'7.69"
$ws.Range("E19").Value = "  -4.89%  "
$ws.Range("E20").Value = "  -1.19%  "
$ws.Range("D21").Value = "'74.61"
$ws.Range("E21").Value = "  -1.74%  "
$ws.Range("E22").Value = "  -6.03%  "
$ws.Range("D23").Value = "'260.22"
$ws.Range("E23").Value = "  -2.72%  "
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").Value = "'9.35"
$ws.Range("E25").Value = "  -7.02%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("D27").Value = "'11.01"
$ws.Range("E27").Value = "  -3.70%  "
$ws.Range("E28").Value = "  +3.08%  "
$ws.Range("D29").Value = "'22.78"
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("D30").Value = "'35.87"
$ws.Range("E30").Value = "  +2.00%  "
$ws.Range("D31").Value = "'0.0905"
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("D32").Value = "'163.63"
$ws.Range("E32").Value = "  -6.60%  "
$ws.Range("D33").Value = "'2.94"
$ws.Range("E33").Value = "  -5.30%  "
$ws.Range("E34").Value = "  -2.77%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.119"
$ws.Range("E35").Value = "  +12.77%  "
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").Value = "'0.130"
$ws.Range("E36").Value = "  -2.07%  "
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("E38").Value = "  -1.09%  "
$ws.Range("E39").Value = "  -5.58%  "
$ws.Range("D40").Value = "'3.60"
$ws.Range("E40").Value = "  -4.82%  "
$ws.Range("D41").Value = "'72.18"
$ws.Range("E41").Value = "  +2.85%  "
$ws.Range("D42").Value = "'98.38"
$ws.Range("E42").Value = "  +7.71%  "
$ws.Range("E43").Value = "  -2.23%  "
$ws.Range("E44").Value = "  -2.85%  "
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "'12.35"
$ws.Range("E46").Value = "  +5.05%  "
$ws.Range("D47").Value = "'112.57"
$ws.Range("E47").Value = "  -5.18%  "
$ws.Range("E48").Value = "  -1.15%  "
$ws.Range("D50").Value = "'74.63"
$ws.Range("E50").Value = "  +3.91%  "
$ws.Range("E51").Value = "  +0.00%  "
